# daily auto push: 2026-01-10 04:03 UTC
# Insert a new data row (2026/01/10, 土, 12, 201) right before the current
# row 605, shifting all subsequent rows down by one (dimension grows from
# A1:D646 to A1:D647).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 605 (and everything below it) down by one row.
$ws.Rows.Item(605).Insert()

# Fill in the newly inserted row. The leading apostrophe on the date forces
# it to be stored as literal text (matching every other date cell in column
# A) instead of being auto-converted into a real Excel date value.
$ws.Range("A605").Value = "'2026/01/10"
$ws.Range("B605").Value = "土"
$ws.Range("C605").Value = 12
$ws.Range("D605").Value = 201

# Drop any formatting Excel may have applied while inserting/typing (e.g. a
# date number format or the "treat as text" quote-prefix flag) so the new
# row's cells stay styled the same as their neighbors.
$ws.Range("A605:D605").ClearFormats()
